# Reorders the comma-separated "Recorded By" names in column G.
# - "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# - "system, backup@backdoor.com, System" -> "backup@backdoor.com, system, System"
# All other values (e.g. "backup@backdoor.com, System", "admin@admin.com, System",
# single-name entries, etc.) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    } elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
}
